$wb = $excel.ActiveWorkbook

# ---- Sheet "snapshot": insert new row 27 (Воронин Кирилл, Torpedo) ----
$ws1 = $wb.Worksheets.Item("snapshot")
$ws1.Rows.Item(27).Insert()

# New row 27 data
$ws1.Range("A27").Value = 'ТОР'
$ws1.Range("B27").Value = 'Торпедо'
$ws1.Range("C27").Value = 'torpedo'
$ws1.Range("D27").Value = 'Воронин Кирилл'
$ws1.Range("E27").Value = '41'
$ws1.Range("F27").Value = 'нападающий'
$ws1.Range("G27").Value = '17354'
$ws1.Range("H27").Value = '1369_ТОР_воронинкирилл'
$ws1.Range("I27").Value = 'injured_active'
$ws1.Range("J27").Value = 'https://www.khl.ru/clubs/torpedo/team/'
$ws1.Range("K27").Value = '2025-12-07T03:02:23.686468+00:00'

# Refresh scraped_at (column K) timestamps for rows 2-26 (unchanged rows)
$ws1.Range("K2").Value = '2025-12-07T03:01:38.169773+00:00'
$ws1.Range("K3").Value = '2025-12-07T03:01:38.169807+00:00'
$ws1.Range("K4").Value = '2025-12-07T03:01:38.169826+00:00'
$ws1.Range("K5").Value = '2025-12-07T03:01:41.117289+00:00'
$ws1.Range("K6").Value = '2025-12-07T03:01:41.117318+00:00'
$ws1.Range("K7").Value = '2025-12-07T03:01:43.416007+00:00'
$ws1.Range("K8").Value = '2025-12-07T03:01:45.710001+00:00'
$ws1.Range("K9").Value = '2025-12-07T03:01:48.150475+00:00'
$ws1.Range("K10").Value = '2025-12-07T03:01:50.932916+00:00'
$ws1.Range("K11").Value = '2025-12-07T03:01:56.322648+00:00'
$ws1.Range("K12").Value = '2025-12-07T03:01:56.322678+00:00'
$ws1.Range("K13").Value = '2025-12-07T03:01:58.633715+00:00'
$ws1.Range("K14").Value = '2025-12-07T03:02:01.428843+00:00'
$ws1.Range("K15").Value = '2025-12-07T03:02:04.125495+00:00'
$ws1.Range("K16").Value = '2025-12-07T03:02:06.505392+00:00'
$ws1.Range("K17").Value = '2025-12-07T03:02:06.505420+00:00'
$ws1.Range("K18").Value = '2025-12-07T03:02:09.254657+00:00'
$ws1.Range("K19").Value = '2025-12-07T03:02:09.254688+00:00'
$ws1.Range("K20").Value = '2025-12-07T03:02:09.254704+00:00'
$ws1.Range("K21").Value = '2025-12-07T03:02:11.669217+00:00'
$ws1.Range("K22").Value = '2025-12-07T03:02:11.669245+00:00'
$ws1.Range("K23").Value = '2025-12-07T03:02:14.499750+00:00'
$ws1.Range("K24").Value = '2025-12-07T03:02:14.499782+00:00'
$ws1.Range("K25").Value = '2025-12-07T03:02:14.499801+00:00'
$ws1.Range("K26").Value = '2025-12-07T03:02:17.251825+00:00'

# Refresh scraped_at (column K) timestamps for rows 28-35 (shifted rows)
$ws1.Range("K28").Value = '2025-12-07T03:02:23.686501+00:00'
$ws1.Range("K29").Value = '2025-12-07T03:02:23.686526+00:00'
$ws1.Range("K30").Value = '2025-12-07T03:02:26.430624+00:00'
$ws1.Range("K31").Value = '2025-12-07T03:02:26.430653+00:00'
$ws1.Range("K32").Value = '2025-12-07T03:02:28.762166+00:00'
$ws1.Range("K33").Value = '2025-12-07T03:02:28.762195+00:00'
$ws1.Range("K34").Value = '2025-12-07T03:02:31.567102+00:00'
$ws1.Range("K35").Value = '2025-12-07T03:02:31.567131+00:00'

# ---- Sheet "returned": clear out the returned-player row (nobody returned this run) ----
$ws2 = $wb.Worksheets.Item("returned")
$ws2.Rows.Item(2).Delete()

# ---- Sheet "new_injured": record the newly injured player ----
$ws3 = $wb.Worksheets.Item("new_injured")
$ws3.Range("A2").Value = 'ТОР'
$ws3.Range("B2").Value = 'Торпедо'
$ws3.Range("C2").Value = 'Воронин Кирилл'
$ws3.Range("D2").Value = '1369_ТОР_воронинкирилл'
$ws3.Range("E2").Value = 'INJURED_NEW'
$ws3.Range("F2").Value = '2025-12-07T11:02:32.076694+08:00'
$ws3.Range("G2").Value = '2025-12-07'
